# "Faaet glemt password knappen til at virke efter hensigten":
# Bjarke's stored password was stale/wrong ("Password1!"); update it to his
# real password so the "forgot password" lookup returns the right value.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("C2").Value = "Hejmeddig1!"

# Touching the page setup (as Excel itself does on a normal save) brings the
# worksheet's header/footer block into the saved XML.
$ws.PageSetup.CenterHeader = ""

# Mirror Excel's full-precision recalculation-on-save behaviour for this
# workbook.
$wb.PrecisionAsDisplayed = $true
$excel.CalculateFullRebuild()
